$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Config")

# Update the test data values (leave balance 1-10 policies)
$ws.Range("D2").Value = "Automation2"
$ws.Range("B2").Value = 42

# Update the active selection to C2
$ws.Range("C2").Select()
